$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.789908885955811
$ws.Range("B1").Value = 4.238265514373779
$ws.Range("C1").Value = 1.90139627456665
$ws.Range("D1").Value = 0.8824335336685181
$ws.Range("E1").Value = 0.4760893881320953
